$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-6 (C2:C6) all become 10908
$ws.Range("C2:C6").Value = 10908

# Rows 7-252 (C7:C252) all become 10820
$ws.Range("C7:C252").Value = 10820
